$wb = $excel.ActiveWorkbook

# Sheets affected: "展览" (sheet1) and "全部类型" (sheet4)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# --- Update "展览" sheet ---
$wsExhibit.Range("F2").Value = 20
$wsExhibit.Range("F3").Value = 54
$wsExhibit.Range("F5").Value = 34
$wsExhibit.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202402/BqSjvdLT1708224342995.png"
$wsExhibit.Range("F7").Value = 2631
$wsExhibit.Range("F8").Value = 1147
$wsExhibit.Range("F9").Value = 232
$wsExhibit.Range("F10").Value = 94
$wsExhibit.Range("F11").Value = 5913
$wsExhibit.Range("F13").Value = 232
$wsExhibit.Range("F14").Value = 581
$wsExhibit.Range("F15").Value = 11585
$wsExhibit.Range("F16").Value = 11788
$wsExhibit.Range("F17").Value = 24
$wsExhibit.Range("F19").Value = 16

# --- Update "全部类型" sheet ---
$wsAll.Range("F2").Value = 20
$wsAll.Range("F3").Value = 54
$wsAll.Range("F5").Value = 34
$wsAll.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202402/BqSjvdLT1708224342995.png"
$wsAll.Range("F7").Value = 2631
$wsAll.Range("F9").Value = 1147
$wsAll.Range("F10").Value = 232
$wsAll.Range("F11").Value = 94
$wsAll.Range("F12").Value = 5913
$wsAll.Range("F14").Value = 232
$wsAll.Range("F15").Value = 581
$wsAll.Range("F16").Value = 11585
$wsAll.Range("F17").Value = 11788
$wsAll.Range("F18").Value = 24
$wsAll.Range("F20").Value = 16
